# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing columns and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, centered, bordered) from G1 onto
# H1, then set its text - Copy() first so the subsequent Value2 write
# isn't clobbered by the paste.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value2 = "Save"

# New "Save" column values, one per data row (2-8).
$saveValues = @(0, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value2 = $saveValues[$i]
}

Write-Output "done"
